$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update country order for Bosnia y Herzegovina / Libano (rows 77 and 78)
$ws.Range("A77").Value = "Bosnia y Herzegovina"
$ws.Range("B77").Value = 23138
$ws.Range("C77").Value = 304
$ws.Range("D77").Value = 15922
$ws.Range("E77").Value = 6526
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 690

$ws.Range("A78").Value = "Libano"
$ws.Range("B78").Value = 22983
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 7312
$ws.Range("E78").Value = 15442
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 229

# Update country order for Namibia / Guayana Francesa (rows 98 and 99)
$ws.Range("A98").Value = "Namibia"
$ws.Range("B98").Value = 9604
$ws.Range("C98").Value = 167
$ws.Range("D98").Value = 5811
$ws.Range("E98").Value = 3695
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 98

$ws.Range("A99").Value = "Guayana Francesa"
$ws.Range("B99").Value = 9494
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 9078
$ws.Range("E99").Value = 353
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 63

# Other numeric updates (rows 4, 19, 24, 52, 65, 92, 104)
$ws.Range("B4").Value = 6641046
$ws.Range("C4").Value = 4799
$ws.Range("E4").Value = 2524738
$ws.Range("G4").Value = 77
$ws.Range("H4").Value = 197498

$ws.Range("B19").Value = 325050
$ws.Range("C19").Value = 643
$ws.Range("D19").Value = 301836
$ws.Range("E19").Value = 18974
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 4240

$ws.Range("B24").Value = 260149
$ws.Range("C24").Value = 424
$ws.Range("E24").Value = 15875
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 9424

$ws.Range("B52").Value = 63310
$ws.Range("C52").Value = 497
$ws.Range("D52").Value = 43894
$ws.Range("E52").Value = 17556
$ws.Range("G52").Value = 5
$ws.Range("H52").Value = 1860

$ws.Range("B65").Value = 42714
$ws.Range("C65").Value = 531
$ws.Range("E65").Value = 11160
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 1117

$ws.Range("B92").Value = 12050
$ws.Range("C92").Value = 47
$ws.Range("E92").Value = 1414

$ws.Range("B104").Value = 8478
$ws.Range("C104").Value = 21
$ws.Range("E104").Value = 2139
$ws.Range("G104").Value = 3
$ws.Range("H104").Value = 219

# Update timestamp string
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 16:54"
